$wb = $excel.ActiveWorkbook

$wsCsv = $wb.Worksheets.Item("forScore CSV Export")
$wsCalc = $wb.Worksheets.Item("IntermediateCalcs")
$wsOut = $wb.Worksheets.Item("Metadata-Output")

# --- Metadata-Output sheet: wrap lookups in IF(ISBLANK(...),"",...) so a blank
# source cell yields "" instead of a literal 0 (missing-terminator / blank
# handling fix referenced in the commit message). ---
$wsOut.Range("A2").Formula = "=IF(ISBLANK('forScore CSV Export'!A2),"""",'forScore CSV Export'!A2)"
$wsOut.Range("B2").Formula = "=IF(ISBLANK('forScore CSV Export'!B2),"""",'forScore CSV Export'!B2)"
$wsOut.Range("C2").Formula = "=IF(ISBLANK('forScore CSV Export'!E2),"""",'forScore CSV Export'!E2)"
$wsOut.Range("D2").Formula = "=IF(ISBLANK('forScore CSV Export'!F2),"""",'forScore CSV Export'!F2)"
$wsOut.Range("E2").Formula = "=IF(ISBLANK('forScore CSV Export'!G2),"""",'forScore CSV Export'!G2)"
$wsOut.Range("F2").Formula = "=IF(ISBLANK('forScore CSV Export'!H2),"""",'forScore CSV Export'!H2)"
$wsOut.Range("G2").Formula = "=IF(ISBLANK('forScore CSV Export'!I2),"""",'forScore CSV Export'!I2)"
$wsOut.Range("H2").Formula = "=IF(IntermediateCalcs!I2<=0,"""",IntermediateCalcs!I2)"

# --- Selection / cursor position updates recorded in the saved views ---
$null = $wsOut.Range("A2").Select()
$null = $wsCalc.Range("H2").Select()

# Restore the originally active sheet/tab ("forScore CSV Export") as the
# last-touched sheet so the saved workbook view still opens there.
$null = $wsCsv.Activate()
